# feat: add 2022-Q1 data
#
# - repurpose the existing "总计" sheet (sheetId 3 / rId3) into the new
#   "2022-Q1" per-fund holdings sheet
# - duplicate it first (preserving sheetPr / pageMargins / styles) so the
#   duplicate becomes the new "总计" summary sheet (sheetId 4 / rId4)
# - fill in the new fund row on "2022-Q1"
# - prepend a "2022-Q1" row to "总计" and push the old rows down

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# Duplicate "总计" right after itself - the copy inherits sheetPr,
# pageMargins, column widths and cell styles exactly.
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item("总计 (2)")

# The original sheet becomes the new quarterly holdings sheet.
$total.Name = "2022-Q1"
# The duplicate becomes the refreshed "总计" sheet.
$newTotal.Name = "总计"

$q1 = $total

# ---------------------------------------------------------------------
# "2022-Q1" sheet: switch the header row from the "总计"-style layout to
# the per-fund layout used by the other quarter sheets, then write the
# single fund holding row.
# ---------------------------------------------------------------------

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

# New header cells E1:H1 need the same bold/centered/bordered style as
# the existing header cells - copy formatting from B1 rather than
# re-declaring it by hand.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2: the fund holding. A2 already carries the index-column style
# (s="2") and value 0 from the old sheet, so it's left untouched.
# The fund code has a leading zero, so it must be forced to text too
# (otherwise it gets coerced into the number 4397).
$q1.Range("B2").Value = "'004397"
$q1.Range("B2").ClearFormats()
$q1.Range("C2").Value = "长盛信息安全量化策略灵活配置混合"

# D2/E2/F2/G2 must stay plain text (e.g. "4.21"), not be coerced into
# numbers. Force text entry with a leading quote, then strip the
# resulting quote-prefix style so the cell matches the unstyled target.
$q1.Range("D2").Value = "'4.21"
$q1.Range("D2").ClearFormats()
$q1.Range("E2").Value = "'29.75"
$q1.Range("E2").ClearFormats()
$q1.Range("F2").Value = "'0.81"
$q1.Range("F2").ClearFormats()
$q1.Range("G2").Value = "'0.0341"
$q1.Range("G2").ClearFormats()

# H2 is a real number.
$q1.Range("H2").Value = 4

# Drop the old row 3 (former "2020-Q4" total row) entirely.
$q1.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# "总计" sheet: insert the new "2022-Q1" summary row above the existing
# rows, pushing "2021-Q1" / "2020-Q4" down by one row.
# ---------------------------------------------------------------------

# Give row 4 the same index-column style as row 3 before filling it in.
$newTotal.Range("A3").Copy()
$newTotal.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Old row 3 ("2020-Q4") data moves to row 4.
$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2020-Q4"
$newTotal.Range("C4").Value = 2
$newTotal.Range("D4").Value = 0.18

# Old row 2 ("2021-Q1") data moves to row 3.
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q1"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.15

# New row 2: "2022-Q1".
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 1
$newTotal.Range("D2").Value = 0.03

# Restore the active sheet/selection similar to the original workbook.
$wb.Worksheets.Item("2020-Q4").Activate()
